$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at G:H, pushing existing G:Q to I:S
$ws.Range("G1:H1").EntireColumn.Insert()

# New header cells
$ws.Range("G1").Value = "MaxIndivGrp"
$ws.Range("H1").Value = "MaxTeamGrp"

# Match the column widths Excel assigned on insert (closest achievable
# given the host's internal 1/6-character snapping grid)
$ws.Columns("G").ColumnWidth = 11.0
$ws.Columns("H").ColumnWidth = 11.5

# Update selection to match the target state (select the two new columns)
$ws.Columns("G:H").Select()
